$d = $word.ActiveDocument

# --- Update the Oakland County "forms you will receive" bullet list ---
#
# Before:
#   * Exhibits to Personal Protection Order (if you uploaded any documents)
#   * Protected Personal Identifying Information forms, if needed
#   * Notice of Hearing on Petition for Personal Protection Order
#   * Personal Protection Order
#   {%p if has_next_friend_petitioning %}
#   * Request for Next Friend and Order
#   {%p endif %}
#
# After:
#   * Exhibits to Personal Protection Order (if you uploaded any documents)
#   * Personal Protection Order
#   * Contact Information Sheet
#   {%p if has_next_friend_petitioning %}
#   * Request for Next Friend and Order
#   {%p endif %}

# Locate the "Notice of Hearing..." bullet paragraph before editing anything,
# so we can unambiguously find & remove the old "Personal Protection Order"
# bullet that immediately follows it (the list has two different bullets with
# that text once the first rename below happens, so find this pairing first).
$noticePara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Notice of Hearing on Petition for Personal Protection Order*") {
        $noticePara = $p
        break
    }
}

if ($noticePara -ne $null) {
    $dupPersonalOrderPara = $noticePara.Next()
    $dupPersonalOrderPara.Range.Delete()
}

# Rename "Protected Personal Identifying Information forms, if needed" -> "Personal Protection Order"
$d.Content.Find.Execute("Protected Personal Identifying Information forms, if needed", $true, $false, $false, $false, $false, $true, 1, $false, "Personal Protection Order", 2)

# Rename "Notice of Hearing on Petition for Personal Protection Order" -> "Contact Information Sheet"
$d.Content.Find.Execute("Notice of Hearing on Petition for Personal Protection Order", $true, $false, $false, $false, $false, $true, 1, $false, "Contact Information Sheet", 2)
